$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Refactored Runmode column: flip the suite rows from YES to NO,
# and the Profile row from YES to Yes (lowercase)
$ws.Range("C2").Value = "NO"
$ws.Range("C3").Value = "NO"
$ws.Range("C4").Value = "NO"
$ws.Range("C5").Value = "NO"
$ws.Range("C6").Value = "NO"
$ws.Range("C7").Value = "NO"
$ws.Range("C8").Value = "Yes"

# Update the active selection on the sheet
$ws.Range("A13").Select()
